# The workbook tracks a Ligand-Receptor (Ccl22-Ackr2) pair across sending
# clusters. The update (new TPM run) changes the existing "Resolving-Mac"
# row's numbers and re-labels it "Inflammatory-Mac", then appends a brand
# new "Resolving-Mac" row with its own freshly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: relabel sending cluster + refresh its numeric columns ---
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Ccl22"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.119110666666667
$ws.Range("H2").Value = 3.357332
$ws.Range("I2").Value = 0.5726510027906514
$ws.Range("J2").Value = 0.5726510027906513
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.062395
$ws.Range("N2").Value = 0.187185
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.06982691004666666
$ws.Range("R2").Value = 0.6284421904199999
$ws.Range("S2").Value = 0.5726510027906514
$ws.Range("T2").Value = 0.5726510027906513

# --- Row 3 (new): the original "Resolving-Mac" row, with new values ---
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Ccl22"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8351523333333333
$ws.Range("H3").Value = 2.505457
$ws.Range("I3").Value = 0.4273489972093487
$ws.Range("J3").Value = 0.4273489972093487
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.062395
$ws.Range("N3").Value = 0.187185
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.05210932983833333
$ws.Range("R3").Value = 0.468983968545
$ws.Range("S3").Value = 0.4273489972093487
$ws.Range("T3").Value = 0.4273489972093487

Write-Host "Edit applied"
